$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-21 Monday" "2024-10-22 Tuesday"

Replace-Text "805×7=" "855×6="
Replace-Text "799×5=" "545×6="
Replace-Text "290×2=" "214×8="
Replace-Text "625×5=" "350×4="
Replace-Text "926×7=" "163×7="

Replace-Text "240×8=" "483×9="
Replace-Text "439×2=" "361×7="
Replace-Text "841×2=" "274×6="
Replace-Text "464×6=" "625×7="
Replace-Text "719×8=" "309×3="

Replace-Text "708×9=" "988×2="
Replace-Text "535×4=" "710×9="
Replace-Text "914×5=" "502×8="
Replace-Text "938×9=" "809×6="
Replace-Text "276×4=" "900×2="

Replace-Text "933×7=" "808×3="
Replace-Text "270×3=" "892×7="
Replace-Text "533×4=" "780×4="
Replace-Text "181×9=" "310×5="
Replace-Text "629×6=" "272×5="

Replace-Text "122×6=" "937×8="
Replace-Text "718×9=" "681×9="
Replace-Text "917×9=" "636×6="
Replace-Text "505×9=" "930×5="
Replace-Text "431×6=" "572×4="
